$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 1684.4
$ws.Cells.Item(17, 10).Value = 1749.3334
$ws.Cells.Item(17, 12).Value = 5248.0002
$ws.Cells.Item(17, 14).Value = -5584.0002
# Row 43
$ws.Cells.Item(43, 8).Value = 500
$ws.Cells.Item(43, 9).Value = 500
$ws.Cells.Item(43, 11).Value = 500
$ws.Cells.Item(43, 13).Value = -431
# Row 58
$ws.Cells.Item(58, 8).Value = 7546.0713
$ws.Cells.Item(58, 9).Value = 724.1667
$ws.Cells.Item(58, 10).Value = 12662.5
$ws.Cells.Item(58, 11).Value = 2172.5001
$ws.Cells.Item(58, 12).Value = 37987.5
$ws.Cells.Item(58, 13).Value = -2022.5001
$ws.Cells.Item(58, 14).Value = -38287.5
# Row 64
$ws.Cells.Item(64, 8).Value = 4581.3335
$ws.Cells.Item(64, 9).Value = 3995.2
$ws.Cells.Item(64, 11).Value = 3995.2
$ws.Cells.Item(64, 13).Value = -3747.2
# Row 67
$ws.Cells.Item(67, 8).Value = 4581.3335
$ws.Cells.Item(67, 9).Value = 3995.2
$ws.Cells.Item(67, 11).Value = 3995.2
$ws.Cells.Item(67, 13).Value = -3137.2
# Row 105
$ws.Cells.Item(105, 8).Value = 51668
$ws.Cells.Item(105, 10).Value = 56001.8
$ws.Cells.Item(105, 12).Value = 56001.8
$ws.Cells.Item(105, 14).Value = -62989.8
# Row 113
$ws.Cells.Item(113, 8).Value = 8472.454
$ws.Cells.Item(113, 9).Value = 8229.125
$ws.Cells.Item(113, 11).Value = 8229.125
$ws.Cells.Item(113, 13).Value = -4975.125
# Row 137
$ws.Cells.Item(137, 8).Value = 1120245.9
$ws.Cells.Item(137, 9).Value = 2005199
$ws.Cells.Item(137, 11).Value = 6015597
$ws.Cells.Item(137, 13).Value = -6013047
# Row 138
$ws.Cells.Item(138, 8).Value = 3933.14
$ws.Cells.Item(138, 9).Value = 2636.375
$ws.Cells.Item(138, 10).Value = 4045.902
$ws.Cells.Item(138, 11).Value = 7909.125
$ws.Cells.Item(138, 12).Value = 12137.706
$ws.Cells.Item(138, 13).Value = -2769.125
$ws.Cells.Item(138, 14).Value = -22417.706
# Row 141
$ws.Cells.Item(141, 8).Value = 8345.5
$ws.Cells.Item(141, 9).Value = 8272.777
$ws.Cells.Item(141, 11).Value = 24818.331
$ws.Cells.Item(141, 13).Value = -19638.331

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 40489.605
$ws.Cells.Item(32, 9).Value = 35170.74
$ws.Cells.Item(32, 10).Value = 56800.8
$ws.Cells.Item(32, 11).Value = 35170.74
$ws.Cells.Item(32, 12).Value = 56800.8
$ws.Cells.Item(32, 13).Value = -34883.74
$ws.Cells.Item(32, 14).Value = -57374.8
# Row 61
$ws.Cells.Item(61, 8).Value = 7147744.5
$ws.Cells.Item(61, 9).Value = 3697.8215
$ws.Cells.Item(61, 11).Value = 3697.8215
$ws.Cells.Item(61, 13).Value = -3485.8215
# Row 63
$ws.Cells.Item(63, 8).Value = 9998.799999999999
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 9998.799999999999
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 9998.799999999999
$ws.Cells.Item(63, 13).Value = ""
$ws.Cells.Item(63, 14).Value = -11370.8
# Row 66
$ws.Cells.Item(66, 8).Value = 9998.799999999999
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 9998.799999999999
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 49994
$ws.Cells.Item(66, 13).Value = ""
$ws.Cells.Item(66, 14).Value = -56858
# Row 74
$ws.Cells.Item(74, 8).Value = 7754.5293
$ws.Cells.Item(74, 9).Value = 2545.16
$ws.Cells.Item(74, 10).Value = 22225
$ws.Cells.Item(74, 11).Value = 2545.16
$ws.Cells.Item(74, 12).Value = 22225
$ws.Cells.Item(74, 13).Value = -1671.16
$ws.Cells.Item(74, 14).Value = -23973
# Row 77
$ws.Cells.Item(77, 8).Value = 7754.5293
$ws.Cells.Item(77, 9).Value = 2545.16
$ws.Cells.Item(77, 10).Value = 22225
$ws.Cells.Item(77, 11).Value = 12725.8
$ws.Cells.Item(77, 12).Value = 111125
$ws.Cells.Item(77, 13).Value = -8357.799999999999
$ws.Cells.Item(77, 14).Value = -119861
# Row 136
$ws.Cells.Item(136, 8).Value = 7147744.5
$ws.Cells.Item(136, 9).Value = 3697.8215
$ws.Cells.Item(136, 11).Value = 11093.4645
$ws.Cells.Item(136, 13).Value = -8543.4645

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 1762
$ws.Cells.Item(22, 9).Value = 1426.8
$ws.Cells.Item(22, 11).Value = 1426.8
$ws.Cells.Item(22, 13).Value = -1253.8
# Row 86
$ws.Cells.Item(86, 8).Value = 2998
$ws.Cells.Item(86, 9).Value = 2729.8333
$ws.Cells.Item(86, 11).Value = 2729.8333
$ws.Cells.Item(86, 13).Value = -1606.8333
# Row 89
$ws.Cells.Item(89, 8).Value = 2998
$ws.Cells.Item(89, 9).Value = 2729.8333
$ws.Cells.Item(89, 11).Value = 13649.1665
$ws.Cells.Item(89, 13).Value = -8033.166499999999
# Row 105
$ws.Cells.Item(105, 8).Value = 1868.625
$ws.Cells.Item(105, 10).Value = 1200
$ws.Cells.Item(105, 12).Value = 1200
$ws.Cells.Item(105, 14).Value = -4694
# Row 134
$ws.Cells.Item(134, 8).Value = 4632605
$ws.Cells.Item(134, 9).Value = 2007.0714
$ws.Cells.Item(134, 11).Value = 6021.2142
$ws.Cells.Item(134, 13).Value = -3486.2142

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 5733.375
$ws.Cells.Item(31, 9).Value = 3487.2
$ws.Cells.Item(31, 10).Value = 7715.294
$ws.Cells.Item(31, 11).Value = 3487.2
$ws.Cells.Item(31, 12).Value = 7715.294
$ws.Cells.Item(31, 13).Value = -3192.2
$ws.Cells.Item(31, 14).Value = -8305.294
# Row 34
$ws.Cells.Item(34, 8).Value = 5733.375
$ws.Cells.Item(34, 9).Value = 3487.2
$ws.Cells.Item(34, 10).Value = 7715.294
$ws.Cells.Item(34, 11).Value = 3487.2
$ws.Cells.Item(34, 12).Value = 7715.294
$ws.Cells.Item(34, 13).Value = -3285.2
$ws.Cells.Item(34, 14).Value = -8119.294
# Row 44
$ws.Cells.Item(44, 8).Value = 20000
$ws.Cells.Item(44, 9).Value = 20000
$ws.Cells.Item(44, 11).Value = 20000
$ws.Cells.Item(44, 13).Value = -19558
# Row 58
$ws.Cells.Item(58, 8).Value = 2500.375
$ws.Cells.Item(58, 9).Value = 1530.1
$ws.Cells.Item(58, 11).Value = 1530.1
$ws.Cells.Item(58, 13).Value = -1327.1
# Row 86
$ws.Cells.Item(86, 8).Value = 45014.04
$ws.Cells.Item(86, 9).Value = 4217.375
$ws.Cells.Item(86, 11).Value = 4217.375
$ws.Cells.Item(86, 13).Value = -3094.375
# Row 89
$ws.Cells.Item(89, 8).Value = 45014.04
$ws.Cells.Item(89, 9).Value = 4217.375
$ws.Cells.Item(89, 11).Value = 21086.875
$ws.Cells.Item(89, 13).Value = -15470.875
# Row 119
$ws.Cells.Item(119, 8).Value = 90000
$ws.Cells.Item(119, 10).Value = 90000
$ws.Cells.Item(119, 12).Value = 90000
$ws.Cells.Item(119, 14).Value = -99676
# Row 136
$ws.Cells.Item(136, 8).Value = 2500.375
$ws.Cells.Item(136, 9).Value = 1530.1
$ws.Cells.Item(136, 11).Value = 4590.299999999999
$ws.Cells.Item(136, 13).Value = -2040.299999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Cells.Item(114, 8).Value = 4138.778
$ws.Cells.Item(114, 10).Value = 4003.5
$ws.Cells.Item(114, 12).Value = 12010.5
$ws.Cells.Item(114, 14).Value = -18518.5
# Row 131
$ws.Cells.Item(131, 8).Value = 14189.5625
$ws.Cells.Item(131, 9).Value = 11240
$ws.Cells.Item(131, 10).Value = 15530.272
$ws.Cells.Item(131, 11).Value = 33720
$ws.Cells.Item(131, 12).Value = 46590.81600000001
$ws.Cells.Item(131, 13).Value = -28680
$ws.Cells.Item(131, 14).Value = -56670.81600000001
# Row 132
$ws.Cells.Item(132, 8).Value = 2782.182
$ws.Cells.Item(132, 10).Value = 1600.6666
$ws.Cells.Item(132, 12).Value = 14405.9994
$ws.Cells.Item(132, 14).Value = -19465.9994

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).Value = ""

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 1850169.9
$ws.Cells.Item(7, 9).Value = 27001.334
$ws.Cells.Item(7, 10).Value = 3673338.2
$ws.Cells.Item(7, 11).Value = 27001.334
$ws.Cells.Item(7, 12).Value = 3673338.2
$ws.Cells.Item(7, 13).Value = -26889.334
$ws.Cells.Item(7, 14).Value = -3673562.2
# Row 46
$ws.Cells.Item(46, 8).Value = 2859.9062
$ws.Cells.Item(46, 9).Value = 2135.15
$ws.Cells.Item(46, 11).Value = 2135.15
$ws.Cells.Item(46, 13).Value = -1947.15
# Row 94
$ws.Cells.Item(94, 8).Value = 42664.5
$ws.Cells.Item(94, 10).Value = 42664.5
$ws.Cells.Item(94, 12).Value = 42664.5
$ws.Cells.Item(94, 14).Value = -44016.5
# Row 109
$ws.Cells.Item(109, 8).Value = 107500
$ws.Cells.Item(109, 10).Value = 107500
$ws.Cells.Item(109, 12).Value = 107500
$ws.Cells.Item(109, 14).Value = -110274
# Row 121
$ws.Cells.Item(121, 8).Value = 112848.56
$ws.Cells.Item(121, 10).Value = 112848.56
$ws.Cells.Item(121, 12).Value = 112848.56
$ws.Cells.Item(121, 14).Value = -116342.56
# Row 122
$ws.Cells.Item(122, 8).Value = 5179.364
$ws.Cells.Item(122, 9).Value = 4630
$ws.Cells.Item(122, 10).Value = 5680.9565
$ws.Cells.Item(122, 11).Value = 13890
$ws.Cells.Item(122, 12).Value = 17042.8695
$ws.Cells.Item(122, 13).Value = -11440
$ws.Cells.Item(122, 14).Value = -21942.8695
# Row 126
$ws.Cells.Item(126, 8).Value = 1850169.9
$ws.Cells.Item(126, 9).Value = 27001.334
$ws.Cells.Item(126, 10).Value = 3673338.2
$ws.Cells.Item(126, 11).Value = 81004.00199999999
$ws.Cells.Item(126, 12).Value = 11020014.6
$ws.Cells.Item(126, 13).Value = -78534.00199999999
$ws.Cells.Item(126, 14).Value = -11024954.6
# Row 136
$ws.Cells.Item(136, 8).Value = 64252.832
$ws.Cells.Item(136, 9).Value = 15340.866
$ws.Cells.Item(136, 10).Value = 145772.78
$ws.Cells.Item(136, 11).Value = 46022.598
$ws.Cells.Item(136, 12).Value = 437318.34
$ws.Cells.Item(136, 13).Value = -43472.598
$ws.Cells.Item(136, 14).Value = -442418.34

$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Cells.Item(27, 8).Value = 79879
$ws.Cells.Item(27, 10).Value = 79879
$ws.Cells.Item(27, 12).Value = 79879
$ws.Cells.Item(27, 14).Value = -80017
# Row 121
$ws.Cells.Item(121, 8).Value = 79829.664
$ws.Cells.Item(121, 10).Value = 79829.664
$ws.Cells.Item(121, 12).Value = 79829.664
$ws.Cells.Item(121, 14).Value = -83323.664
# Row 126
$ws.Cells.Item(126, 8).Value = 1702.9412
$ws.Cells.Item(126, 9).Value = 1589.2858
$ws.Cells.Item(126, 11).Value = 4767.857400000001
$ws.Cells.Item(126, 13).Value = -2297.857400000001
# Row 136
$ws.Cells.Item(136, 8).Value = 4636524
$ws.Cells.Item(136, 9).Value = 5001.6
$ws.Cells.Item(136, 11).Value = 15004.8
$ws.Cells.Item(136, 13).Value = -12454.8
